# Generate Report for handback
# This script updates the localization-status workbook to reflect that the
# zh-cn and de-de handoffs have now been handed back and are in sync with
# en-US: the Status is updated, and the "Latest Target File" / "Latest
# Handback File" / "Latest Handback DateTime" columns are populated for the
# two real source-file rows on each language sheet.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$mdFile      = "78795ea4-5631-49d0-875f-f4a01c40f99e.md"
$xlfZhCn     = "78795ea4-5631-49d0-875f-f4a01c40f99e.1d95289e63adc3a3b376d22e4b035eb255bb0bbe.zh-cn.xlf"
$xlfDeDe     = "78795ea4-5631-49d0-875f-f4a01c40f99e.1d95289e63adc3a3b376d22e4b035eb255bb0bbe.de-de.xlf"

$targetUrlZhCn   = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/64d831f6722a70a5956e42bba5c8cf145baf4234/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$mdFile"
$handbackUrlZhCn = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/64d831f6722a70a5956e42bba5c8cf145baf4234/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$xlfZhCn"
$targetUrlDeDe   = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/64d831f6722a70a5956e42bba5c8cf145baf4234/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$mdFile"
$handbackUrlDeDe = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/64d831f6722a70a5956e42bba5c8cf145baf4234/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$xlfDeDe"

$handbackDtZhCn = "2016-01-18 07:02:09"
$handbackDtDeDe = "2016-01-18 07:02:26"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# 1. Update the "Status" column wherever it shows "Ready for handoff" -> now
#    handed back and in sync with en-US (rows 2 and 3 on every sheet).
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn.Range("B2").Value = $newStatus
$wsZhCn.Range("B3").Value = $newStatus

$wsDeDe.Range("B2").Value = $newStatus
$wsDeDe.Range("B3").Value = $newStatus

# 2. Populate "Latest Target File" (E) and "Latest Handback File" (F) for the
#    two real source rows (2 and 3) on the zh-cn sheet, plus the actual
#    handback datetime (G).
$wsZhCn.Range("E2").Value = $mdFile
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E2"), $targetUrlZhCn, "", "", $mdFile) | Out-Null

$wsZhCn.Range("F2").Value = $xlfZhCn
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), $handbackUrlZhCn, "", "", $xlfZhCn) | Out-Null

$wsZhCn.Range("G2").Value = $handbackDtZhCn

$wsZhCn.Range("E3").Value = $mdFile
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E3"), $targetUrlZhCn, "", "", $mdFile) | Out-Null

$wsZhCn.Range("F3").Value = $xlfZhCn
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), $handbackUrlZhCn, "", "", $xlfZhCn) | Out-Null

$wsZhCn.Range("G3").Value = $handbackDtZhCn

# 3. Same for the de-de sheet.
$wsDeDe.Range("E2").Value = $mdFile
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E2"), $targetUrlDeDe, "", "", $mdFile) | Out-Null

$wsDeDe.Range("F2").Value = $xlfDeDe
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), $handbackUrlDeDe, "", "", $xlfDeDe) | Out-Null

$wsDeDe.Range("G2").Value = $handbackDtDeDe

$wsDeDe.Range("E3").Value = $mdFile
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E3"), $targetUrlDeDe, "", "", $mdFile) | Out-Null

$wsDeDe.Range("F3").Value = $xlfDeDe
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), $handbackUrlDeDe, "", "", $xlfDeDe) | Out-Null

$wsDeDe.Range("G3").Value = $handbackDtDeDe
